$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the FilesTab query text in cell B4 (drop the `File Type` and
#     `Breed` columns from the Neo4j query, per the corrected ICDC Breed
#     1-14 scripts commit). ---
$newFilesTabQuery = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Cavalier King Charles Spaniel']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS ``File Name``,
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@
$ws.Range("B4").Value = $newFilesTabQuery.TrimEnd("`r", "`n")

# Row 4 wraps the (now shorter) query text, so its auto-fit height shrinks.
$ws.Rows.Item(4).RowHeight = 217.5

# --- Window/view changes: scrolled to A4 (best effort - engine only
#     persists zoom + selection), zoomed to 70%, selection moved to B4. ---
$win = $excel.ActiveWindow
$win.Zoom = 70
[void]$ws.Range("B4").Select()
